$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header text updates: issue number 48 -> 49, week dates 11/28-12/4 -> 12/5-12/11
$ws.Range("C1").Value = "Volume 29   Number  49"
$ws.Range("C8").Value = "Report Covering the Week  12/5/2022  Through  12/11/2022"

# Cells whose type flips between number and text this week (style must follow)
$ws.Range("D22").Copy($ws.Range("C22"))
$ws.Range("C14").Copy($ws.Range("D23"))
$ws.Range("E22").Copy($ws.Range("E23"))
$ws.Range("C16").Copy($ws.Range("C23"))
$ws.Range("C23").Value = 1
$ws.Range("C16").Copy($ws.Range("C27"))
$ws.Range("C27").Value = 1
$ws.Range("C16").Copy($ws.Range("C28"))
$ws.Range("C28").Value = 2
$ws.Range("C16").Copy($ws.Range("C29"))
$ws.Range("C29").Value = 1
$ws.Range("D16").Copy($ws.Range("D28"))
$ws.Range("D28").Value = 1
$ws.Range("D16").Copy($ws.Range("D29"))
$ws.Range("D29").Value = 1
$ws.Range("E16").Copy($ws.Range("E28"))
$ws.Range("E28").Value = 100
$ws.Range("E16").Copy($ws.Range("E29"))
$ws.Range("E29").Value = 0

# Plain value refreshes (style unchanged)
$ws.Range("G15").Value = 3
$ws.Range("H15").Value = -66.666666666666
$ws.Range("L15").Value = 76.190476190476
$ws.Range("M15").Value = 8.823529411764
$ws.Range("N15").Value = -39.344262295082
$ws.Range("C16").Value = 7
$ws.Range("D16").Value = 4
$ws.Range("E16").Value = 75
$ws.Range("F16").Value = 26
$ws.Range("G16").Value = 21
$ws.Range("H16").Value = 23.809523809523
$ws.Range("I16").Value = 309
$ws.Range("J16").Value = 234
$ws.Range("K16").Value = 32.051282051282
$ws.Range("L16").Value = 30.379746835443
$ws.Range("M16").Value = -27.294117647058
$ws.Range("N16").Value = -78.777472527472
$ws.Range("C17").Value = 9
$ws.Range("D17").Value = 8
$ws.Range("E17").Value = 12.5
$ws.Range("F17").Value = 29
$ws.Range("G17").Value = 24
$ws.Range("H17").Value = 20.833333333333
$ws.Range("I17").Value = 359
$ws.Range("J17").Value = 306
$ws.Range("K17").Value = 17.320261437908
$ws.Range("L17").Value = 19.666666666666
$ws.Range("M17").Value = -0.277777777777
$ws.Range("N17").Value = -55.348258706467
$ws.Range("C18").Value = 8
$ws.Range("D18").Value = 8
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 22
$ws.Range("G18").Value = 22
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 287
$ws.Range("J18").Value = 254
$ws.Range("K18").Value = 12.992125984252
$ws.Range("L18").Value = -11.419753086419
$ws.Range("M18").Value = -35.214446952595
$ws.Range("N18").Value = -77.700077700077
$ws.Range("C19").Value = 8
$ws.Range("D19").Value = 14
$ws.Range("E19").Value = -42.857142857142
$ws.Range("F19").Value = 43
$ws.Range("G19").Value = 50
$ws.Range("H19").Value = -14
$ws.Range("I19").Value = 666
$ws.Range("J19").Value = 536
$ws.Range("K19").Value = 24.253731343283
$ws.Range("L19").Value = 29.06976744186
$ws.Range("M19").Value = 129.655172413793
$ws.Range("N19").Value = 24.022346368715
$ws.Range("C20").Value = 6
$ws.Range("E20").Value = 100
$ws.Range("F20").Value = 16
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 234
$ws.Range("J20").Value = 166
$ws.Range("K20").Value = 40.963855421686
$ws.Range("L20").Value = 41.818181818181
$ws.Range("M20").Value = 43.558282208589
$ws.Range("N20").Value = -74.200661521499
$ws.Range("C21").Value = 38
$ws.Range("D21").Value = 37
$ws.Range("E21").Value = 2.702702702702
$ws.Range("G21").Value = 136
$ws.Range("H21").Value = 0.735294117647
$ws.Range("I21").Value = 1895
$ws.Range("J21").Value = 1530
$ws.Range("K21").Value = 23.856209150326
$ws.Range("L21").Value = 20.854591836734
$ws.Range("M21").Value = 9.664351851851
$ws.Range("N21").Value = -62.784760408483
$ws.Range("F22").Value = 3
$ws.Range("G22").Value = 1
$ws.Range("H22").Value = 200
$ws.Range("M22").Value = -38.461538461538
$ws.Range("F23").Value = 3
$ws.Range("G23").Value = 3
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 28
$ws.Range("K23").Value = -15.151515151515
$ws.Range("L23").Value = -9.677419354838
$ws.Range("M23").Value = 21.739130434782
$ws.Range("D24").Value = 20
$ws.Range("E24").Value = 15
$ws.Range("F24").Value = 75
$ws.Range("G24").Value = 96
$ws.Range("H24").Value = -21.875
$ws.Range("I24").Value = 1004
$ws.Range("J24").Value = 920
$ws.Range("K24").Value = 9.130434782608
$ws.Range("L24").Value = -9.46798917944
$ws.Range("M24").Value = 38.674033149171
$ws.Range("C25").Value = 18
$ws.Range("D25").Value = 19
$ws.Range("E25").Value = -5.263157894736
$ws.Range("F25").Value = 39
$ws.Range("G25").Value = 49
$ws.Range("H25").Value = -20.408163265306
$ws.Range("I25").Value = 558
$ws.Range("J25").Value = 503
$ws.Range("K25").Value = 10.934393638171
$ws.Range("L25").Value = 19.742489270386
$ws.Range("M25").Value = -25.500667556742
$ws.Range("G26").Value = 5
$ws.Range("H26").Value = -80
$ws.Range("L26").Value = 54.545454545454
$ws.Range("E27").Value = 0
$ws.Range("G27").Value = 5
$ws.Range("H27").Value = -40
$ws.Range("I27").Value = 63
$ws.Range("J27").Value = 49
$ws.Range("K27").Value = 28.571428571428
$ws.Range("L27").Value = 6.779661016949
$ws.Range("F28").Value = 3
$ws.Range("H28").Value = 200
$ws.Range("I28").Value = 17
$ws.Range("J28").Value = 38
$ws.Range("K28").Value = -55.263157894736
$ws.Range("L28").Value = 41.666666666666
$ws.Range("M28").Value = -60.465116279069
$ws.Range("N28").Value = -92.307692307692
$ws.Range("F29").Value = 2
$ws.Range("H29").Value = 100
$ws.Range("I29").Value = 16
$ws.Range("J29").Value = 28
$ws.Range("K29").Value = -42.857142857142
$ws.Range("L29").Value = 33.333333333333
$ws.Range("M29").Value = -51.515151515151
$ws.Range("N29").Value = -92.156862745098
